# Seminario: ajustando numeracao paginas slides
# 1) Refresh the cached "today" text of the datetimeFigureOut fields
#    (slide master + every slide layout) from 30/05/2024 to 03/06/2024.
# 2) Fix the page-number rectangles on slides 3..27: each slide's
#    duplicated/stale page number is bumped so it matches the slide's
#    real position in the deck (slide N -> "N").

$p = $ppt.ActivePresentation

# --- 1) Date placeholder on the slide master ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "30/05/2024") {
            $tr.Text = "03/06/2024"
        }
    }
}

# --- Date placeholder on every slide layout ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "30/05/2024") {
                $tr.Text = "03/06/2024"
            }
        }
    }
}

# --- Date placeholder on the notes master (only reachable through the
#     HeadersFooters/DateAndTime object in this environment) ---
$notesMaster = $p.NotesMaster
$nmDate = $notesMaster.HeadersFooters.DateAndTime
$nmDate.Text = "03/06/2024"

# --- 2) Page-number rectangles on slides 3..27 ---
for ($si = 3; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -match '^[0-9]+$') {
                $old = [int]$t
                # Slide 26 carries a leftover duplicate shape showing "1"
                # (from a copy/paste) that is not the visible page number
                # and must stay untouched.
                if ($si -eq 26 -and $old -eq 1) {
                    continue
                }
                $shp.TextFrame.TextRange.Text = [string]$si
            }
        }
    }
}
